# Updates the cryptocurrency price/volume table (rows 2-51) with
# refreshed values for the "Price" (D) and "Volume(1h)" (E) columns,
# as produced by the scheduled GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.295.19"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "3.342.34"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.69"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.79"
$ws.Range("E6").Value = "  +2.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("E9").Value = "  +5.26%  "
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.18"
$ws.Range("E11").Value = "  +6.52%  "
$ws.Range("E12").Value = "  +2.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "692.97"
$ws.Range("E13").Value = "  +3.84%  "
$ws.Range("D14").Value = "3.883.55"
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.44"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "68.295.15"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("E17").Value = "  +1.41%  "
$ws.Range("D18").Value = "3.344.10"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("E20").Value = "  +3.03%  "
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.46"
$ws.Range("E22").Value = "  +1.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.00"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "100.12"
$ws.Range("E24").Value = "  +2.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.93"
$ws.Range("E25").Value = "  +2.68%  "
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.53"
$ws.Range("E27").Value = "  +3.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.08"
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("E29").Value = "  +1.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.99"
$ws.Range("E30").Value = "  -4.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "568.13"
$ws.Range("E31").Value = "  -2.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.05"
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("E33").Value = "  +2.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "57.52"
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "3.706.36"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.33"
$ws.Range("E37").Value = "  +2.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.137"
$ws.Range("E38").Value = "  +4.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.90"
$ws.Range("E39").Value = "  +7.03%  "
$ws.Range("E40").Value = "  +3.59%  "
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.337"
$ws.Range("E42").Value = "  +1.90%  "
$ws.Range("E43").Value = "  +2.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.26"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0413"
$ws.Range("E45").Value = "  +2.23%  "
$ws.Range("E46").Value = "  +3.70%  "
$ws.Range("E47").Value = "  +1.23%  "
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "130.93"
$ws.Range("E50").Value = "  +3.23%  "
$ws.Range("E51").Value = "  +0.65%  "
